$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 1.988074333333333
$ws.Cells.Item(2, 8).Value = 5.964223
$ws.Cells.Item(2, 9).Value = 0.01657769708907969
$ws.Cells.Item(2, 10).Value = 0.01657769708907968
$ws.Cells.Item(2, 13).Value = 29.47402433333333
$ws.Cells.Item(2, 14).Value = 88.422073
$ws.Cells.Item(2, 15).Value = 0.295877356230023
$ws.Cells.Item(2, 16).Value = 0.295877356230023
$ws.Cells.Item(2, 17).Value = 58.5965512771421
$ws.Cells.Item(2, 18).Value = 527.368961494279
$ws.Cells.Item(2, 19).Value = 0.004904965187099045
$ws.Cells.Item(2, 20).Value = 0.004904965187099044

# Row 3
$ws.Cells.Item(3, 7).Value = 1.988074333333333
$ws.Cells.Item(3, 8).Value = 5.964223
$ws.Cells.Item(3, 9).Value = 0.01657769708907969
$ws.Cells.Item(3, 10).Value = 0.01657769708907968
$ws.Cells.Item(3, 15).Value = 0.1818061388681701
$ws.Cells.Item(3, 16).Value = 0.1818061388681701
$ws.Cells.Item(3, 17).Value = 36.00550199051344
$ws.Cells.Item(3, 18).Value = 324.049517914621
$ws.Cells.Item(3, 19).Value = 0.003013927099091681
$ws.Cells.Item(3, 20).Value = 0.003013927099091681

# Row 4
$ws.Cells.Item(4, 7).Value = 1.988074333333333
$ws.Cells.Item(4, 8).Value = 5.964223
$ws.Cells.Item(4, 9).Value = 0.01657769708907969
$ws.Cells.Item(4, 10).Value = 0.01657769708907968
$ws.Cells.Item(4, 13).Value = 7.238098333333333
$ws.Cells.Item(4, 14).Value = 21.714295
$ws.Cells.Item(4, 15).Value = 0.07266023040422054
$ws.Cells.Item(4, 16).Value = 0.07266023040422054
$ws.Cells.Item(4, 17).Value = 14.38987751864278
$ws.Cells.Item(4, 18).Value = 129.508897667785
$ws.Cells.Item(4, 19).Value = 0.001204539290063906
$ws.Cells.Item(4, 20).Value = 0.001204539290063906

# Row 5
$ws.Cells.Item(5, 7).Value = 1.988074333333333
$ws.Cells.Item(5, 8).Value = 5.964223
$ws.Cells.Item(5, 9).Value = 0.01657769708907969
$ws.Cells.Item(5, 10).Value = 0.01657769708907968
$ws.Cells.Item(5, 13).Value = 44.79281599999999
$ws.Cells.Item(5, 14).Value = 134.378448
$ws.Cells.Item(5, 15).Value = 0.4496562744975863
$ws.Cells.Item(5, 16).Value = 0.4496562744975864
$ws.Cells.Item(5, 17).Value = 89.05144780732266
$ws.Cells.Item(5, 18).Value = 801.4630302659039
$ws.Cells.Item(5, 19).Value = 0.007454265512825054
$ws.Cells.Item(5, 20).Value = 0.007454265512825053

# Row 6
$ws.Cells.Item(6, 9).Value = 0.7746030815641455
$ws.Cells.Item(6, 10).Value = 0.7746030815641454
$ws.Cells.Item(6, 13).Value = 29.47402433333333
$ws.Cells.Item(6, 14).Value = 88.422073
$ws.Cells.Item(6, 15).Value = 0.295877356230023
$ws.Cells.Item(6, 16).Value = 0.295877356230023
$ws.Cells.Item(6, 17).Value = 2737.95985922587
$ws.Cells.Item(6, 18).Value = 24641.63873303283
$ws.Cells.Item(6, 19).Value = 0.2291875119008282
$ws.Cells.Item(6, 20).Value = 0.2291875119008282

# Row 7
$ws.Cells.Item(7, 9).Value = 0.7746030815641455
$ws.Cells.Item(7, 10).Value = 0.7746030815641454
$ws.Cells.Item(7, 15).Value = 0.1818061388681701
$ws.Cells.Item(7, 16).Value = 0.1818061388681701
$ws.Cells.Item(7, 19).Value = 0.1408275954145636
$ws.Cells.Item(7, 20).Value = 0.1408275954145635

# Row 8
$ws.Cells.Item(8, 9).Value = 0.7746030815641455
$ws.Cells.Item(8, 10).Value = 0.7746030815641454
$ws.Cells.Item(8, 13).Value = 7.238098333333333
$ws.Cells.Item(8, 14).Value = 21.714295
$ws.Cells.Item(8, 15).Value = 0.07266023040422054
$ws.Cells.Item(8, 16).Value = 0.07266023040422054
$ws.Cells.Item(8, 17).Value = 672.3758679734756
$ws.Cells.Item(8, 18).Value = 6051.38281176128
$ws.Cells.Item(8, 19).Value = 0.05628283837827004
$ws.Cells.Item(8, 20).Value = 0.05628283837827004

# Row 9
$ws.Cells.Item(9, 9).Value = 0.7746030815641455
$ws.Cells.Item(9, 10).Value = 0.7746030815641454
$ws.Cells.Item(9, 13).Value = 44.79281599999999
$ws.Cells.Item(9, 14).Value = 134.378448
$ws.Cells.Item(9, 15).Value = 0.4496562744975863
$ws.Cells.Item(9, 16).Value = 0.4496562744975864
$ws.Cells.Item(9, 17).Value = 4160.983610608982
$ws.Cells.Item(9, 18).Value = 37448.85249548083
$ws.Cells.Item(9, 19).Value = 0.3483051358704837
$ws.Cells.Item(9, 20).Value = 0.3483051358704836

# Row 10
$ws.Cells.Item(10, 7).Value = 23.741365
$ws.Cells.Item(10, 8).Value = 71.22409500000001
$ws.Cells.Item(10, 9).Value = 0.1979690350870239
$ws.Cells.Item(10, 10).Value = 0.1979690350870239
$ws.Cells.Item(10, 13).Value = 29.47402433333333
$ws.Cells.Item(10, 14).Value = 88.422073
$ws.Cells.Item(10, 15).Value = 0.295877356230023
$ws.Cells.Item(10, 16).Value = 0.295877356230023
$ws.Cells.Item(10, 17).Value = 699.7535697165483
$ws.Cells.Item(10, 18).Value = 6297.782127448935
$ws.Cells.Item(10, 19).Value = 0.0585745547169573
$ws.Cells.Item(10, 20).Value = 0.0585745547169573

# Row 11
$ws.Cells.Item(11, 7).Value = 23.741365
$ws.Cells.Item(11, 8).Value = 71.22409500000001
$ws.Cells.Item(11, 9).Value = 0.1979690350870239
$ws.Cells.Item(11, 10).Value = 0.1979690350870239
$ws.Cells.Item(11, 15).Value = 0.1818061388681701
$ws.Cells.Item(11, 16).Value = 0.1818061388681701
$ws.Cells.Item(11, 17).Value = 429.9737441566184
$ws.Cells.Item(11, 18).Value = 3869.763697409565
$ws.Cells.Item(11, 19).Value = 0.03599198588462912
$ws.Cells.Item(11, 20).Value = 0.03599198588462912

# Row 12
$ws.Cells.Item(12, 7).Value = 23.741365
$ws.Cells.Item(12, 8).Value = 71.22409500000001
$ws.Cells.Item(12, 9).Value = 0.1979690350870239
$ws.Cells.Item(12, 10).Value = 0.1979690350870239
$ws.Cells.Item(12, 13).Value = 7.238098333333333
$ws.Cells.Item(12, 14).Value = 21.714295
$ws.Cells.Item(12, 15).Value = 0.07266023040422054
$ws.Cells.Item(12, 16).Value = 0.07266023040422054
$ws.Cells.Item(12, 17).Value = 171.8423344375583
$ws.Cells.Item(12, 18).Value = 1546.581009938025
$ws.Cells.Item(12, 19).Value = 0.01438447570232438
$ws.Cells.Item(12, 20).Value = 0.01438447570232438

# Row 13
$ws.Cells.Item(13, 7).Value = 23.741365
$ws.Cells.Item(13, 8).Value = 71.22409500000001
$ws.Cells.Item(13, 9).Value = 0.1979690350870239
$ws.Cells.Item(13, 10).Value = 0.1979690350870239
$ws.Cells.Item(13, 13).Value = 44.79281599999999
$ws.Cells.Item(13, 14).Value = 134.378448
$ws.Cells.Item(13, 15).Value = 0.4496562744975863
$ws.Cells.Item(13, 16).Value = 0.4496562744975864
$ws.Cells.Item(13, 17).Value = 1063.44259403384
$ws.Cells.Item(13, 18).Value = 9570.98334630456
$ws.Cells.Item(13, 19).Value = 0.08901801878311313
$ws.Cells.Item(13, 20).Value = 0.08901801878311315

# Row 14
$ws.Cells.Item(14, 7).Value = 1.301204666666667
$ws.Cells.Item(14, 8).Value = 3.903614
$ws.Cells.Item(14, 9).Value = 0.01085018625975097
$ws.Cells.Item(14, 10).Value = 0.01085018625975097
$ws.Cells.Item(14, 13).Value = 29.47402433333333
$ws.Cells.Item(14, 14).Value = 88.422073
$ws.Cells.Item(14, 15).Value = 0.295877356230023
$ws.Cells.Item(14, 16).Value = 0.295877356230023
$ws.Cells.Item(14, 17).Value = 38.35173800798022
$ws.Cells.Item(14, 18).Value = 345.165642071822
$ws.Cells.Item(14, 19).Value = 0.003210324425138439
$ws.Cells.Item(14, 20).Value = 0.003210324425138438

# Row 15
$ws.Cells.Item(15, 7).Value = 1.301204666666667
$ws.Cells.Item(15, 8).Value = 3.903614
$ws.Cells.Item(15, 9).Value = 0.01085018625975097
$ws.Cells.Item(15, 10).Value = 0.01085018625975097
$ws.Cells.Item(15, 15).Value = 0.1818061388681701
$ws.Cells.Item(15, 16).Value = 0.1818061388681701
$ws.Cells.Item(15, 17).Value = 23.56578244093089
$ws.Cells.Item(15, 18).Value = 212.092041968378
$ws.Cells.Item(15, 19).Value = 0.001972630469885796
$ws.Cells.Item(15, 20).Value = 0.001972630469885796

# Row 16
$ws.Cells.Item(16, 7).Value = 1.301204666666667
$ws.Cells.Item(16, 8).Value = 3.903614
$ws.Cells.Item(16, 9).Value = 0.01085018625975097
$ws.Cells.Item(16, 10).Value = 0.01085018625975097
$ws.Cells.Item(16, 13).Value = 7.238098333333333
$ws.Cells.Item(16, 14).Value = 21.714295
$ws.Cells.Item(16, 15).Value = 0.07266023040422054
$ws.Cells.Item(16, 16).Value = 0.07266023040422054
$ws.Cells.Item(16, 17).Value = 9.418247329125556
$ws.Cells.Item(16, 18).Value = 84.76422596213
$ws.Cells.Item(16, 19).Value = 0.0007883770335622135
$ws.Cells.Item(16, 20).Value = 0.0007883770335622134

# Row 17
$ws.Cells.Item(17, 7).Value = 1.301204666666667
$ws.Cells.Item(17, 8).Value = 3.903614
$ws.Cells.Item(17, 9).Value = 0.01085018625975097
$ws.Cells.Item(17, 10).Value = 0.01085018625975097
$ws.Cells.Item(17, 13).Value = 44.79281599999999
$ws.Cells.Item(17, 14).Value = 134.378448
$ws.Cells.Item(17, 15).Value = 0.4496562744975863
$ws.Cells.Item(17, 16).Value = 0.4496562744975864
$ws.Cells.Item(17, 17).Value = 58.28462121234134
$ws.Cells.Item(17, 18).Value = 524.561590911072
$ws.Cells.Item(17, 19).Value = 0.004878854331164522
$ws.Cells.Item(17, 20).Value = 0.004878854331164522
